$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated "Price" values are plain decimals (e.g. "0.999", "315.05").
# The Price column stores text (European-style thousands separators like
# "43.110.00" appear elsewhere in the same column), so force NumberFormat to
# Text on the cells that would otherwise be auto-parsed as numbers before
# assigning their values, keeping them literal text like the rest of the column.
foreach ($addr in @("D4","D5","D6","D10","D12","D17","D19","D20","D22","D23","D26","D29","D30","D31","D32","D37","D39","D42","D47","D49","D50")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.110.00"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.563.23"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "315.05"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").Value = "96.64"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "35.47"
$ws.Range("E10").Value = "  -2.38%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "7.43"
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("D13").Value = "2.957.28"
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "2.587.14"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "0.841"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "43.112.01"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("D19").Value = "6.81"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").Value = "12.56"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").Value = "0.0₃0960"
$ws.Range("E21").Value = "  -1.13%  "
$ws.Range("D22").Value = "69.19"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").Value = "253.43"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("E24").Value = "  -0.33%  "
$ws.Range("D26").Value = "26.81"
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "40.05"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "10.17"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").Value = "5.82"
$ws.Range("E31").Value = "  -4.42%  "
$ws.Range("D32").Value = "155.10"
$ws.Range("E32").Value = "  +0.78%  "
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("D37").Value = "18.94"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("D39").Value = "2.45"
$ws.Range("E39").Value = "  +6.56%  "
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("E41").Value = "  -6.28%  "
$ws.Range("D42").Value = "3.97"
$ws.Range("E42").Value = "  +3.91%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("E45").Value = "  -2.74%  "
$ws.Range("D46").Value = "2.003.71"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "8.90"
$ws.Range("E47").Value = "  +0.83%  "
$ws.Range("D48").Value = "2.810.84"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("D49").Value = "82.58"
$ws.Range("E49").Value = "  -3.54%  "
$ws.Range("D50").Value = "74.63"
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  +1.81%  "
